$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "States testing with JDBC"
$ws.Range("B3").Value = "FAILED"
$ws.Range("C3").Value = "chrome"

$ws.Range("A4").Value = "States testing with JDBC"
$ws.Range("B4").Value = "FAILED"
$ws.Range("C4").Value = "chrome"

$ws.Range("A5").Value = "States testing with JDBC"
$ws.Range("B5").Value = "FAILED"
$ws.Range("C5").Value = "chrome"

$ws.Range("A6").Value = "States testing with JDBC"
$ws.Range("B6").Value = "FAILED"
$ws.Range("C6").Value = "chrome"
